# Apply "x" signal markers to column G for the specified rows in the
# current_quarter_classification sheet (test fixed portfolio 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("current_quarter_classification")

$rows = @(11,26,37,43,63,74,75,79,81,87,95,97,106,124,132,134,143,144,151,155,157,161,170,171,175,183,188,194,202,213,235,239,243,246,251,269,273,282,288,290,292,297,298,309,310,321,328)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "x"
}
